$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Row 16: new worker ANDRES FELIPE SANCHEZ MARTINEZ (moved up from row 19)
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "1143384090"
$ws.Range("D16").Value = "ANDRES FELIPE SANCHEZ MARTINEZ"
$ws.Range("E16").Value = "1907"
$ws.Range("F16").Value = 22083
$ws.Range("G16").Value = 828116

# Row 17: GLADIS MARIA GALEANO POLO, periodo 1810
$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1041258536"
$ws.Range("D17").Value = "GLADIS MARIA GALEANO POLO"
$ws.Range("E17").Value = "1810"
$ws.Range("F17").Value = 31249
$ws.Range("G17").Value = 877803

# Row 18: GLADIS MARIA GALEANO POLO, periodo 1809
$ws.Range("B18").Value = "CC"
$ws.Range("C18").Value = "1041258536"
$ws.Range("D18").Value = "GLADIS MARIA GALEANO POLO"
$ws.Range("E18").Value = "1809"
$ws.Range("F18").Value = 31249
$ws.Range("G18").Value = 877803

# Row 19: GLADIS MARIA GALEANO POLO, periodo 1808
$ws.Range("B19").Value = "CC"
$ws.Range("C19").Value = "1041258536"
$ws.Range("D19").Value = "GLADIS MARIA GALEANO POLO"
$ws.Range("E19").Value = "1808"
$ws.Range("F19").Value = 15625
$ws.Range("G19").Value = 877803
